$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

$data = @(
  @(2, 0.03140023387760406, -0.007160354249104994),
  @(3, 0.03490555789464717, -0.0005695455839876962),
  @(4, 0.03276250820358476, 0.0001228048630723677),
  @(5, 0.03818356007532743, -0.0009018418026918162),
  @(6, 0.0163811357666355, -0.002745069710323067),
  @(7, 0.01660597256470389, -0.01617615620323543),
  @(8, 0.03280510886006088, -0.01158638203317208),
  @(9, 0.03195049235708722, -0.005066666666666553),
  @(10, 0.03199013463464139, 0.007139290140011356),
  @(11, 0.02932700192909973, 0.02339103417665322),
  @(12, 0.01623676687524422, 0.04780992639020476),
  @(13, 0.01680063389776837, -0.01567177319950697),
  @(14, 0.008192342910670922, 0.0002022244691608499),
  @(15, 0.007993776517429468, -0.000710563713879675),
  @(16, 0.0312816620504122, 0.003631576558173366),
  @(17, 0.03203699535676512, -0.003579198540248574),
  @(18, 0.03195475242273484, -0.001155401502022091),
  @(19, 0.03289942198009272, -0.008578550386844075),
  @(20, 0.02633608083900573, -0.0003594617060952432),
  @(21, 0.0306289253250726, 0.0108796439389256),
  @(22, 0.03359061763141771, -0.005136334812936072),
  @(23, 0.03238975245941873, -0.01034664172560928),
  @(24, 0.0168793267770923, 0.01805945036455414),
  @(25, 0.01551018901201267, 0.0404364080262456),
  @(26, 0.0321705957488805, -0.01258000441403673),
  @(27, 0.03224088683206609, -0.01623760341195224),
  @(28, 0.03242111127599143, -0.004858072027943972),
  @(29, 0.03197889279473797, 0.01509769094138536),
  @(30, 0.03310426013664871, -0.009008042895442325),
  @(31, 0.03218763601147094, 0.003562447611064501),
  @(32, 0.03535854487517653, -0.01654618473895586),
  @(33, 0.03066371586119476, -0.01569506726457404),
  @(34, 0.04590753243507482, 0.0002319916483008999),
  @(35, 0.03090369955934355, -0.007030334591850096),
  @(36, 0.03402017425088521, 0.0003895787679570883),
  @(37, 1, -0.001325590427348367)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value = $item[1]
    $ws.Cells.Item($r, 5).Value = $item[2]
}

$ws.Protect()
